$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data one column right and one row down,
# which preserves the original numeric cell values/widths for
# the columns/rows that already existed.
$ws.Columns("A:A").Insert()
$ws.Rows("1:1").Insert()

# New header row (row 1) labels
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# New row-label column (column A) for each data row
$ws.Range("A2").Value = "CyclomaticComplexity(CC) & EffortToImplement"
$ws.Range("A3").Value = "NbOperands & EffortToImplement"
$ws.Range("A4").Value = "NbOperators & EffortToImplement"
$ws.Range("A5").Value = "ProgramLevel & ProgramLevel"
$ws.Range("A6").Value = "EffortToImplement & NbOperands"
$ws.Range("A7").Value = "EffortToImplement & NbOperators"

# Widen the new first column for the long row labels (closest
# achievable value given the column-width rounding granularity).
$ws.Columns("A:A").ColumnWidth = 53.6
